# Update the "想去人数" (number of people wanting to go) counts in column F
# for rows 2-5 on both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 635
    $ws.Range("F3").Value = 3841
    $ws.Range("F4").Value = 107
    $ws.Range("F5").Value = 724
}
